# Red book.xlsx - fix weighted-average tire/ethanol data mix-up.
#
# The data-entry columns H ("Snow Tires (0-no,1-yes)") and I ("Ethanol
# (0-no,1-yes)") were swapped for every fill-up row from row 370 through
# the end of the data (row 702) - those rows had H=1/I=0 when they should
# have read H=0/I=1. Fix the values, then leave the sheet scrolled/
# selected on the corrected Ethanol column and touch Page Setup the way
# the author's Excel session did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the swapped Snow Tires / Ethanol columns for rows 370-702 -----
$ws.Range("H370:H702").Value = 0
$ws.Range("I370:I702").Value = 1

# --- Touch Page Setup (author re-saved with an explicit page setup) ----
$ws.PageSetup.Orientation = 1

# --- Update the view: scroll back to the top and select the corrected
#     Ethanol column so the selection/scroll state matches the saved file
$ws.Application.Goto($ws.Range("A1"), $true) | Out-Null
$ws.Range("H2:H702").Select() | Out-Null
